$wb = $excel.ActiveWorkbook

# --- Metadata sheet: revert Version / Date / Contact values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"

# --- "Include from FSIII" sheet: drop the two UUID concept rows so the
#     remaining rows (C / D / blank / System URI) shift up two places ---
$inc = $wb.Worksheets.Item("Include from FSIII")

# Current layout (rows 1-7):
#   1: Concept           | Description
#   2: <uuid 1>          | (blank)
#   3: <uuid 2>          | (blank)
#   4: C                 | (blank)
#   5: D                 | (blank)
#   6: (blank)           | (blank)
#   7: System URI        | urn:oid:1.2.208.176.2.21
#
# Target layout (rows 1-5), obtained by deleting rows 2 and 3 so the rest
# shift up naturally (preserving their existing cell encodings verbatim):
#   1: Concept           | Description
#   2: C                 | (blank)
#   3: D                 | (blank)
#   4: (blank)           | (blank)
#   5: System URI        | urn:oid:1.2.208.176.2.21

$inc.Rows.Item(3).Delete()
$inc.Rows.Item(2).Delete()
